$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel rows 8-15 (column B) - the shared-string table gets two new
# entries ("line7","line8") inserted right after "line6", which shifts
# the effective label every later row resolves to by two positions.
$ws.Range("B8").Value  = "line7"
$ws.Range("B9").Value  = "line8"
$ws.Range("B10").Value = "extr1"
$ws.Range("B11").Value = "extr2"
$ws.Range("B12").Value = "extr3"
$ws.Range("B13").Value = "extr4"
$ws.Range("B14").Value = "extr5"
$ws.Range("B15").Value = "extr6"

# Update the numeric / boolean contents of rows 8-15
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# Add two new rows (16 and 17) for line7 / line8 contingencies,
# copying the formatting of an existing data row first.
$ws.Range("A2").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A2").Copy($ws.Range("A17"))
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
